$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": update F3 value and move selection to A3 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F3").Value = 5218.3100000000004
$wsSummary.Range("A3").Select()

# --- Sheet "Transactions": update row 2/3 values, fix number formats, move selection to B2 ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 210
$wsTrans.Range("C2").Value = 42024
$wsTrans.Range("E2").NumberFormat = "#,##0.00"
$wsTrans.Range("E2").Value = 4438.6000000000004
$wsTrans.Range("F2").NumberFormat = "#,##0.00"
$wsTrans.Range("F2").Value = 4126.2700000000004
$wsTrans.Range("G2").Value = 312.33
$wsTrans.Range("J2").NumberFormat = "#,##0.00"
$wsTrans.Range("J2").Value = 45873.73
$wsTrans.Range("A3").Value = 207
$wsTrans.Range("B2").Select()

# --- Sheet "Input": move selection to B2 ---
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B2").Select()

# --- Sheet "Repayment Schedule": move selection to B7 ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Range("B7").Select()

# Re-activate the Transactions sheet so it remains the tab-selected / active sheet
$wsTrans.Activate()
$wsTrans.Range("B2").Select()
